$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: replace stale test placeholders with the real record pulled from the database
$ws.Range("A10").Value = "Vesselin"
$ws.Range("B10").Value = "Netzov"
$ws.Range("C10").Value = 51
$ws.Range("D10").Value = "Nokia"
$ws.Range("E10").Value = "Developer"
$ws.Range("F10").Value = 5000

# Row 11: A11 keeps a formula (now pointing at A7), remaining cells become literals
$ws.Range("A11").Formula = "=A7"
$ws.Range("B11").Value = "Netzov"
$ws.Range("C11").Value = 111
$ws.Range("D11").Value = "Don’t know"
$ws.Range("E11").Value = "DDD"
$ws.Range("F11").Value = 111

# Row 7: "Test" -> "Test1"
$ws.Range("A7").Value = "Test1"

# Row 8: mirrors row 7
$ws.Range("A8").Value = "Test1"

# Row 12 is no longer part of the data set
$ws.Rows(12).Delete()

# Restore the user's last selection before save
$ws.Range("E12").Select()
